$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 : (Intercept)
$ws.Range("C2").Value = 534
$ws.Range("D2").Value = 40738.47926089869
$ws.Range("F2").Value = 0.9870615962606853
$ws.Range("G2").Value = 17.46874775196302

# Row 3 : A
$ws.Range("C3").Value = 534
$ws.Range("D3").Value = 5.3635623236671846
$ws.Range("E3").Value = 0.020939571881875785
$ws.Range("F3").Value = 0.009944242989941836
$ws.Range("G3").Value = 0.2004407563508374
$ws.Range("I3").Value = "*"

# Row 4 : B
$ws.Range("C4").Value = 534
$ws.Range("D4").Value = 20.711434969594457
$ws.Range("E4").Value = 0.0000000021724756438246118
$ws.Range("F4").Value = 0.07198683282012498
$ws.Range("G4").Value = 0.5570311211740968

# Row 5 : A:B
$ws.Range("C5").Value = 534
$ws.Range("D5").Value = 0.4704832940980626
$ws.Range("E5").Value = 0.6249589828629885
$ws.Range("F5").Value = 0.0017590101468532553
$ws.Range("G5").Value = 0.08395498116783973
